# Text updates as supplied by PM&C.
# Update the "Description" sheet of the Skills - higher level qualifications
# dashboard metadata workbook:
#   - refresh the preliminary-data footnote from the 2014/July-2016 wording
#     to the 2015/July-2017 wording
#   - split the old combined "Source: NCVER (unpublished) ..." cell into a
#     label cell ("Source") in column A and the source text itself in
#     column B

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

$ws.Range("B8").Value = "The 2015 data are preliminary. Final figures will be released by the National Centre for Vocational Education Research in July 2017.  "

$ws.Range("A9").Value = "Source"
$ws.Range("B9").Value = "NCVER (unpublished) National VET provider collection."

# Row 9 no longer wraps onto two lines now that the source text moved out
# of column A, so its height shrinks back to the single-line row height.
$ws.Rows.Item(9).RowHeight = 13.8

# Reflect the active selection resting on the updated source cell.
$ws.Range("B9").Select() | Out-Null
